# 150 report without template
# The "Sheet4" example sheet (the *-CITYFROM/-CITYTO merge demo) drops its
# "func= with white font" helper note (columns M:N) and instead enables the
# FIRST() aggregation function on the CITYFROM/CITYTO/CARRNAME merge tags so
# the report works correctly without the extra template note.

$wb = $excel.ActiveWorkbook
$prevActive = $wb.ActiveSheet

# Sheet2's cursor was left parked on the old J1:M1 merge; move it off to N1
# (mirrors the same "park the cursor past the used range" tidy-up applied to
# Sheet4 below).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("N1").Select()

$ws = $wb.Worksheets.Item("Sheet4")

# Drop the demonstration note that lived in M1:N1 (and its leftover
# formatting) - no longer needed.
$ws.Range("M1:N1").Clear()

# Switch the merge tags to also run FIRST() so the from/to city (and the
# carrier name) are aggregated per the merged group.
$ws.Range("D2").Value = "{R-T-CARRNAME;merge=X}"
$ws.Range("B3").Value = "{R-T-CITYFROM;func=FIRST;merge=X}"
$ws.Range("C3").Value = "{R-T-CITYTO;func=FIRST;merge=X}"
$ws.Range("D3").Value = "{R-T-CARRNAME;func=FIRST}"

# The old M3/N3 helper formulas for CITYFROM/CITYTO FIRST() are now folded
# into the merge tags above, so clear their old locations.
$ws.Range("M3:N3").Clear()

$ws.Range("N1").Select()

# Restore whichever sheet/tab was active before we touched Sheet4 so the
# workbook's active-tab state is left untouched.
$prevActive.Activate()
